$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rename every occurrence of "Sarah" to "Farah" across the team / actor cells
$ws.Range("D9").Value = "Hugo,Farah,Alvyn"
$ws.Range("D20").Value = "Hugo,Farah,Alvyn"
$ws.Range("D6").Value = "Farah"
$ws.Range("D11").Value = "Farah"
$ws.Range("D14").Value = "Farah"
$ws.Range("D15").Value = "Farah"
$ws.Range("D19").Value = "Farah"
$ws.Range("D2").Value = "Farah Villard, Hugo Poissonnet, Alvyn Silou"

# Update the sheet view's active cell / selection & scroll position
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("E4").Select()
